# Update global battery, wind/solar capacity to EU start year
#
# 1) SYBGSaWC: shift the "start year" references forward by one column
#    (Global Renewables Outlook col C -> D, Wind col B -> C).
# 2) BGSaWC: drop the old first-year column (B, 2021) so the sheet now
#    starts at 2022, shifting every later year left by one column.

$wb = $excel.ActiveWorkbook

$sy = $wb.Worksheets.Item("SYBGSaWC")
$sy.Range("B7").Formula = "='Global Renewables Outlook'!D7*(1-Wind!C4)"
$sy.Range("B8").Formula = "='Global Renewables Outlook'!D6"
$sy.Range("B15").Formula = "='Global Renewables Outlook'!D7*Wind!C4"

$bg = $wb.Worksheets.Item("BGSaWC")
$bg.Columns("B").Delete()

$bg.Activate()
$bg.Range("D22").Select()
